# "Generate Report for Handback"
#
# The ba8ba465-25d9-4545-815a-ec9d8f67e65d file failed its handback
# transform for both the zh-cn and de-de locales. Update the status
# shown on the Overview sheet plus each locale sheet, and record the
# specific mismatch in the "Error Detail" column (K) of row 3 on each
# locale sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the ba8ba465... file; column B is zh-cn
# status, column C is de-de status.
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# zh-cn sheet: row 3 Status + Error Detail.
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K3").Value = "Handback file name: fsan4xmx.xb5 is different with handoff file name: ba8ba465-25d9-4545-815a-ec9d8f67e65d.b77942c259e2062b2fa8121c8a563895d538fe77.zh-cn."

# de-de sheet: row 3 Status + Error Detail.
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K3").Value = "Handback file name: fsan4xmx.xb5 is different with handoff file name: ba8ba465-25d9-4545-815a-ec9d8f67e65d.b77942c259e2062b2fa8121c8a563895d538fe77.de-de."
